$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.730.17'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '3.186.62'
$ws.Range("E3").Value = '  -2.76%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.71'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '615.49'
$ws.Range("E6").Value = '  -2.84%  '
$ws.Range("E7").Value = '  +1.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.689'
$ws.Range("E8").Value = '  -5.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '3.182.86'
$ws.Range("E10").Value = '  -2.77%  '
$ws.Range("E11").Value = '  -1.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.177'
$ws.Range("E12").Value = '  -5.14%  '
$ws.Range("E13").Value = '  -5.10%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.775.22'
$ws.Range("E14").Value = '  -2.63%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '89.518.40'
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.89'
$ws.Range("E16").Value = '  -5.15%  '
$ws.Range("E17").Value = '  -5.50%  '
$ws.Range("D18").Value = '3.170.74'
$ws.Range("E18").Value = '  -2.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.27'
$ws.Range("E19").Value = '  +3.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.41'
$ws.Range("E20").Value = '  -5.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.59'
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000196'
$ws.Range("E22").Value = '  +36.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.60'
$ws.Range("E23").Value = '  -4.66%  '
$ws.Range("E24").Value = '  -6.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.13'
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.71'
$ws.Range("E26").Value = '  -5.70%  '
$ws.Range("D27").Value = '3.344.48'
$ws.Range("E27").Value = '  -3.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '75.36'
$ws.Range("E28").Value = '  -3.03%  '
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("E30").Value = '  -7.02%  '
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.14'
$ws.Range("E32").Value = '  +30.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.43'
$ws.Range("E33").Value = '  -5.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '533.44'
$ws.Range("E34").Value = '  -6.81%  '
$ws.Range("E35").Value = '  -3.38%  '
$ws.Range("E36").Value = '  -6.28%  '
$ws.Range("E37").Value = '  -8.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '21.97'
$ws.Range("E38").Value = '  -4.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.31'
$ws.Range("E39").Value = '  +2.19%  '
$ws.Range("E40").Value = '  -9.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("E43").Value = '  -6.19%  '
$ws.Range("E44").Value = '  -8.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '149.36'
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '172.64'
$ws.Range("E46").Value = '  -4.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.59'
$ws.Range("E47").Value = '  -2.66%  '
$ws.Range("E48").Value = '  -8.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.23'
$ws.Range("E49").Value = '  -8.51%  '
$ws.Range("E50").Value = '  -4.78%  '
$ws.Range("E51").Value = '  -3.98%  '
